$p = $ppt.ActivePresentation

# Slide 5: the "Test your solution in Judge" textbox currently spells the
# hyperlink out across three separate runs (".../Recursive" + "-" +
# "Algorithms"); join them back into a single run with the full URL text,
# keeping the same run formatting/hyperlink.
$slide5 = $p.Slides.Item(5)
$linkBox = $slide5.Shapes.Item("TextBox 10")
$origLinkHeight = $linkBox.Height
$linkRange = $linkBox.TextFrame.TextRange
$fullLinkText = $linkRange.Text
$fullUrl = "https://judge.softuni.bg/Contests/2726/Recursive-Algorithms"
$oldUrlMarker = "https://judge.softuni.bg/Contests/2726/Recursive"
$urlStart = $fullLinkText.IndexOf($oldUrlMarker) + 1
$urlLen = $fullLinkText.Length - $urlStart + 1
$urlRange = $linkRange.Characters($urlStart, $urlLen)
$urlRange.Text = $fullUrl
# merging the runs can nudge the autosized textbox height by a hair; put it
# back exactly where it was
$linkBox.Height = $origLinkHeight

# Slide 7: the title currently reads "...: Дърво на рекурсия" - the word
# should be lower-cased to "дърво" (kept as its own run in the middle of
# the title, matching how PowerPoint splits out a reformatted/respelled
# word).
$slide7 = $p.Slides.Item(7)
$titleShape = $slide7.Shapes.Item("Rectangle 2")
$titleRange = $titleShape.TextFrame.TextRange
$titleText = $titleRange.Text
$oldWord = "Дърво"
$newWord = "дърво"
$wordStart = $titleText.IndexOf($oldWord) + 1
$wordRange = $titleRange.Characters($wordStart, $oldWord.Length)
$wordRange.Text = $newWord
